$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.000003504660980979679

$ws.Range("A3").Value = 0.000003435254257055931
$ws.Range("C3").Value = 4.0677971839904785
$ws.Range("D3").Value = 1.8305089473724365

$ws.Range("A4").Value = 0.00000006940678076716722
$ws.Range("C4").Value = 0.16949200630187988
$ws.Range("D4").Value = 4.237287998199463
